# Add two new worksheets ("Sheet2", "Sheet3") to the workbook, each holding
# a single text value in A1, then restore the focus back to Sheet1 so its
# original tab stays the selected one.

$wb = $excel.ActiveWorkbook

# New sheets get appended after the last existing sheet so the tab order
# ends up Sheet1, Sheet2, Sheet3.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "sheet2"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = "sheet3"

# Keep Sheet1 as the active/selected tab, like in the source workbook.
$wb.Worksheets.Item("Sheet1").Activate()
